$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B7").Value = 2.0
$ws.Range("C7").Value = $false
$ws.Range("B8").Value = 3.0
$ws.Range("C8").Value = $false
$ws.Range("B9").Value = 4.0
$ws.Range("C9").Value = $true

$ws.Range("B6:C6").Copy()
$ws.Range("B7:C9").PasteSpecial(-4122)
